# LHCb data 8 TeV
# Add "obs"/"units" header columns (X,Y) and fill "dsig/dpT"/"pb" values for all data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1): new "obs" and "units" columns
$ws.Range("X1").Value = "obs"
$ws.Range("Y1").Value = "units"

# Data rows 2-15: populate the new columns with constant values
for ($r = 2; $r -le 15; $r++) {
    $ws.Cells.Item($r, 24).Value = "dsig/dpT"
    $ws.Cells.Item($r, 25).Value = "pb"
}

# Update the window's visible/selected range to mirror the authored view state
$excel.ActiveWindow.ScrollColumn = 12
$ws.Range("X1:Y1048576").Select()
